$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16, shifting the existing rows 16-18
# (1000000/10000000/50000000 data points) down to rows 17-19.
$ws.Rows.Item(16).Insert()

# Fill in the newly inserted row 16 with the extra "500000" data point,
# matching the Interpolation/Reduction/Size of the surrounding rows.
$ws.Range("A16").Value = 500000
$ws.Range("B16").Value = $ws.Range("B17").Value2
$ws.Range("C16").Value = $ws.Range("C17").Value2
$ws.Range("D16").Value = $ws.Range("D17").Value2
$ws.Range("E16").Value = 26.207897778505899

# Update the selection to match the saved workbook state.
[void]$ws.Range("E21").Select()
